$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a cell as literal text. Several "Price" values in this sheet
# use "." as a thousands separator (e.g. "30.403.68") or otherwise look like
# numbers Excel would auto-convert (e.g. "1.000", "69.00"); force text via a
# temporary "@" number format, then restore the default style afterwards so
# the cell keeps its original (unstyled) appearance.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.403.68"
$ws.Range("E2").Value = "  -0.44%  "
Set-TextValue $ws.Range("D3") "1.924.18"
$ws.Range("E3").Value = "  +3.74%  "
Set-TextValue $ws.Range("D4") "1.000"
$ws.Range("E4").Value = "  -0.01%  "
Set-TextValue $ws.Range("D5") "240.18"
$ws.Range("E5").Value = "  +2.67%  "
Set-TextValue $ws.Range("D6") "1.000"
$ws.Range("E6").Value = "  +0.01%  "
Set-TextValue $ws.Range("D7") "0.4737"
$ws.Range("E7").Value = "  -0.23%  "
Set-TextValue $ws.Range("D8") "0.2852"
$ws.Range("E8").Value = "  +3.56%  "
Set-TextValue $ws.Range("D9") "0.06601"
$ws.Range("E9").Value = "  +4.46%  "
Set-TextValue $ws.Range("D10") "19.19"
$ws.Range("E10").Value = "  +7.98%  "
Set-TextValue $ws.Range("D11") "105.00"
$ws.Range("E11").Value = "  +24.11%  "
Set-TextValue $ws.Range("D12") "1.912.63"
$ws.Range("E12").Value = "  +2.88%  "
Set-TextValue $ws.Range("D13") "0.07577"
$ws.Range("E13").Value = "  +1.66%  "
Set-TextValue $ws.Range("D14") "5.136"
$ws.Range("E14").Value = "  +2.76%  "
Set-TextValue $ws.Range("D15") "0.6539"
$ws.Range("E15").Value = "  +4.22%  "
Set-TextValue $ws.Range("D16") "297.84"
$ws.Range("E16").Value = "  +21.43%  "
Set-TextValue $ws.Range("D17") "30.400.64"
$ws.Range("E17").Value = "  -0.33%  "
Set-TextValue $ws.Range("D18") "1.001"
$ws.Range("E18").Value = "  +0.12%  "
Set-TextValue $ws.Range("D19") "12.94"
$ws.Range("E19").Value = "  +1.74%  "
$ws.Range("B20").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C20").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
Set-TextValue $ws.Range("D20") "2.173.06"
$ws.Range("E20").Value = "  +3.28%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue $ws.Range("D21") "0.000007524"
$ws.Range("E21").Value = "  +2.50%  "
Set-TextValue $ws.Range("D22") "1.000"
$ws.Range("E22").Value = "  +0.21%  "
Set-TextValue $ws.Range("D23") "5.216"
$ws.Range("E23").Value = "  +5.22%  "
Set-TextValue $ws.Range("D24") "6.291"
$ws.Range("E24").Value = "  +5.93%  "
Set-TextValue $ws.Range("D25") "9.236"
$ws.Range("E25").Value = "  +1.05%  "
Set-TextValue $ws.Range("D26") "166.27"
$ws.Range("E26").Value = "  +2.04%  "
Set-TextValue $ws.Range("D27") "19.57"
$ws.Range("E27").Value = "  +8.72%  "
Set-TextValue $ws.Range("D28") "2.036"
$ws.Range("E28").Value = "  +8.42%  "
Set-TextValue $ws.Range("D29") "0.1116"
$ws.Range("E29").Value = "  +9.23%  "
Set-TextValue $ws.Range("D30") "1.361"
$ws.Range("E30").Value = "  +0.13%  "
Set-TextValue $ws.Range("D31") "4.102"
$ws.Range("E31").Value = "  +2.10%  "
Set-TextValue $ws.Range("D32") "3.926"
$ws.Range("E32").Value = "  +2.34%  "
Set-TextValue $ws.Range("D33") "0.05010"
$ws.Range("E33").Value = "  +3.46%  "
Set-TextValue $ws.Range("D34") "0.7397"
$ws.Range("E34").Value = "  +4.95%  "
$ws.Range("E35").Value = "  +0.27%  "
Set-TextValue $ws.Range("D36") "0.9999"
$ws.Range("E36").Value = "  +0.05%  "
Set-TextValue $ws.Range("D37") "2.717"
$ws.Range("E37").Value = "  +0.83%  "
Set-TextValue $ws.Range("D38") "0.01953"
$ws.Range("E38").Value = "  +2.66%  "
Set-TextValue $ws.Range("D39") "2.693"
$ws.Range("E39").Value = "  +0.28%  "
Set-TextValue $ws.Range("D40") "2.039"
$ws.Range("E40").Value = "  +1.33%  "
Set-TextValue $ws.Range("D41") "0.8721"
$ws.Range("E41").Value = "  -0.60%  "
Set-TextValue $ws.Range("D42") "107.31"
$ws.Range("E42").Value = "  +0.49%  "
Set-TextValue $ws.Range("D43") "5.824"
$ws.Range("E43").Value = "  +4.73%  "
Set-TextValue $ws.Range("D44") "1.000"
$ws.Range("E44").Value = "  +0.00%  "
Set-TextValue $ws.Range("D45") "69.00"
$ws.Range("E45").Value = "  +9.56%  "
Set-TextValue $ws.Range("D46") "0.4124"
$ws.Range("E46").Value = "  +1.43%  "
Set-TextValue $ws.Range("D47") "7.274"
$ws.Range("E47").Value = "  +0.99%  "
Set-TextValue $ws.Range("D48") "9.220"
$ws.Range("E48").Value = "  +7.61%  "
$ws.Range("E49").Value = "  -0.45%  "
Set-TextValue $ws.Range("D50") "34.65"
$ws.Range("E50").Value = "  +2.96%  "
Set-TextValue $ws.Range("D51") "0.05623"
$ws.Range("E51").Value = "  +1.54%  "
